$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts existing rows 4-9 down to 5-10)
$ws.Rows.Item(4).Insert()

# Fill in the new row 4 with the France validation ROI data
# (set D4 "france_valid" before A4's coordinate text so new shared strings
# are appended in the same order as the target workbook)
$ws.Cells.Item(4, 4).Value = "france_valid"
$ws.Cells.Item(4, 1).Value = "[[[3.4341882784815203, 48.48852574424973],`n           [3.4341882784815203, 48.28423633983985],`n           [3.9986108859033953, 48.28423633983985],`n           [3.9986108859033953, 48.48852574424973]]]"
$ws.Cells.Item(4, 2).Value = "2018-09-01"
$ws.Cells.Item(4, 3).Value = "suger beat, kolza, barley"
$ws.Rows.Item(4).RowHeight = 115.2

# Set the new column D width to fit the "france_valid" / country labels
$ws.Columns.Item(4).ColumnWidth = 44.5

# Update the selection/view state
$ws.Range("B4").Select()
